$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 456.2857
$ws.Range("I53").Value = 757.1429000000001
$ws.Range("J53").Value = 155.42857
$ws.Range("K53").Value = 757.1429000000001
$ws.Range("L53").Value = 155.42857
$ws.Range("M53").Value = -120.1429000000001
$ws.Range("N53").Value = -1429.42857
$ws.Range("H62").Value = 100003040
$ws.Range("I62").Value = 111113160
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 111113160
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -111112536
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 100003040
$ws.Range("I65").Value = 111113160
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 555565800
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -555562680
$ws.Range("N65").Value = -66240
$ws.Range("H116").Value = 6030554.5
$ws.Range("I116").Value = 7267078
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 7267078
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -7263636
$ws.Range("N116").Value = -9384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1947.0476
$ws.Range("I61").Value = 1787.4445
$ws.Range("K61").Value = 1787.4445
$ws.Range("M61").Value = -1575.4445
$ws.Range("H97").Value = 1358.8334
$ws.Range("I97").Value = 1472.3334
$ws.Range("J97").Value = 337.33334
$ws.Range("K97").Value = 1472.3334
$ws.Range("L97").Value = 337.33334
$ws.Range("M97").Value = -976.3334
$ws.Range("N97").Value = -1329.33334
$ws.Range("H110").Value = 5808.9165
$ws.Range("I110").Value = 6558.316
$ws.Range("J110").Value = 2961.2
$ws.Range("K110").Value = 6558.316
$ws.Range("L110").Value = 2961.2
$ws.Range("M110").Value = -4513.316
$ws.Range("N110").Value = -7051.2
$ws.Range("H122").Value = 2302.9546
$ws.Range("I122").Value = 2085.8333
$ws.Range("K122").Value = 6257.499899999999
$ws.Range("M122").Value = -3807.499899999999
$ws.Range("H132").Value = 2347.6875
$ws.Range("I132").Value = 2005.3334
$ws.Range("J132").Value = 3374.75
$ws.Range("K132").Value = 6016.0002
$ws.Range("L132").Value = 10124.25
$ws.Range("M132").Value = -3486.0002
$ws.Range("N132").Value = -15184.25
$ws.Range("H136").Value = 1947.0476
$ws.Range("I136").Value = 1787.4445
$ws.Range("K136").Value = 5362.333500000001
$ws.Range("M136").Value = -2812.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 94.111115
$ws.Range("I80").Value = 120
$ws.Range("J80").Value = 90.875
$ws.Range("K80").Value = 120
$ws.Range("L80").Value = 90.875
$ws.Range("M80").Value = 878
$ws.Range("N80").Value = -2086.875
$ws.Range("H83").Value = 94.111115
$ws.Range("I83").Value = 120
$ws.Range("J83").Value = 90.875
$ws.Range("K83").Value = 600
$ws.Range("L83").Value = 454.375
$ws.Range("M83").Value = 4392
$ws.Range("N83").Value = -10438.375
$ws.Range("H94").Value = 6456.943
$ws.Range("I94").Value = 353.9565
$ws.Range("J94").Value = 18154.334
$ws.Range("K94").Value = 353.9565
$ws.Range("L94").Value = 18154.334
$ws.Range("M94").Value = 97.04349999999999
$ws.Range("N94").Value = -19056.334
$ws.Range("H99").Value = 1551
$ws.Range("I99").Value = 1445.5555
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1445.5555
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 52.44450000000006
$ws.Range("N99").Value = -5496
$ws.Range("H105").Value = 3991
$ws.Range("I105").Value = 5624.75
$ws.Range("J105").Value = 2901.8333
$ws.Range("K105").Value = 5624.75
$ws.Range("L105").Value = 2901.8333
$ws.Range("M105").Value = -3877.75
$ws.Range("N105").Value = -6395.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 240390.58
$ws.Range("I86").Value = 557544.9
$ws.Range("J86").Value = 2524.8333
$ws.Range("K86").Value = 557544.9
$ws.Range("L86").Value = 2524.8333
$ws.Range("M86").Value = -556421.9
$ws.Range("N86").Value = -4770.8333
$ws.Range("H89").Value = 240390.58
$ws.Range("I89").Value = 557544.9
$ws.Range("J89").Value = 2524.8333
$ws.Range("K89").Value = 2787724.5
$ws.Range("L89").Value = 12624.1665
$ws.Range("M89").Value = -2782108.5
$ws.Range("N89").Value = -23856.1665
$ws.Range("H105").Value = 743.75
$ws.Range("I105").Value = 658.3333
$ws.Range("K105").Value = 658.3333
$ws.Range("M105").Value = 1088.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 85.25
$ws.Range("I14").Value = 85.25
$ws.Range("K14").Value = 255.75
$ws.Range("M14").Value = -82.75
$ws.Range("H33").Value = 5421.1577
$ws.Range("I33").Value = 432
$ws.Range("J33").Value = 7203
$ws.Range("K33").Value = 2592
$ws.Range("L33").Value = 43218
$ws.Range("M33").Value = -2309
$ws.Range("N33").Value = -43784
$ws.Range("H80").Value = 1171.174
$ws.Range("I80").Value = 998.8
$ws.Range("J80").Value = 1219.0555
$ws.Range("K80").Value = 2996.4
$ws.Range("L80").Value = 3657.1665
$ws.Range("M80").Value = -2060.4
$ws.Range("N80").Value = -5529.166499999999
$ws.Range("H83").Value = 1171.174
$ws.Range("I83").Value = 998.8
$ws.Range("J83").Value = 1219.0555
$ws.Range("K83").Value = 8989.199999999999
$ws.Range("L83").Value = 10971.4995
$ws.Range("M83").Value = -4309.199999999999
$ws.Range("N83").Value = -20331.4995
$ws.Range("H92").Value = 194.6842
$ws.Range("I92").Value = 185.57143
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 556.71429
$ws.Range("L92").Value = 600
$ws.Range("M92").Value = 691.28571
$ws.Range("N92").Value = -3096
$ws.Range("H96").Value = 9733.333000000001
$ws.Range("J96").Value = 9733.333000000001
$ws.Range("L96").Value = 29199.999
$ws.Range("N96").Value = -33317.999
$ws.Range("H98").Value = 372.63635
$ws.Range("I98").Value = 285.7143
$ws.Range("J98").Value = 524.75
$ws.Range("K98").Value = 857.1428999999999
$ws.Range("L98").Value = 1574.25
$ws.Range("M98").Value = 640.8571000000001
$ws.Range("N98").Value = -4570.25
$ws.Range("H107").Value = 264.25
$ws.Range("I107").Value = 296.25
$ws.Range("J107").Value = 248.25
$ws.Range("K107").Value = 888.75
$ws.Range("L107").Value = 744.75
$ws.Range("M107").Value = 1031.25
$ws.Range("N107").Value = -4584.75
$ws.Range("H131").Value = 744.2353000000001
$ws.Range("I131").Value = 360.30768
$ws.Range("J131").Value = 981.9048
$ws.Range("K131").Value = 1080.92304
$ws.Range("L131").Value = 2945.7144
$ws.Range("M131").Value = 3959.07696
$ws.Range("N131").Value = -13025.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2158
$ws.Range("I126").Value = 2245.1428
$ws.Range("J126").Value = 1975
$ws.Range("K126").Value = 6735.428400000001
$ws.Range("L126").Value = 5925
$ws.Range("M126").Value = -4265.428400000001
$ws.Range("N126").Value = -10865
$ws.Range("H132").Value = 3301
$ws.Range("I132").Value = 3075
$ws.Range("J132").Value = 3401.4443
$ws.Range("K132").Value = 9225
$ws.Range("L132").Value = 10204.3329
$ws.Range("M132").Value = -6695
$ws.Range("N132").Value = -15264.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 665.61536
$ws.Range("I107").Value = 528.1111
$ws.Range("J107").Value = 975
$ws.Range("K107").Value = 1584.3333
$ws.Range("L107").Value = 2925
$ws.Range("M107").Value = 335.6667000000002
$ws.Range("N107").Value = -6765
$ws.Range("H109").Value = 8300
$ws.Range("J109").Value = 8300
$ws.Range("L109").Value = 8300
$ws.Range("N109").Value = -11074
